$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Tache (B) and Description supplementaire (E) text for rows 3-25 (reworded to add trailing periods / corrections) ---
$ws.Range("B3").Value = "Lecture du CDC avec les experts"
$ws.Range("E3").Value = "Reçu le cahier des charges + premier entretient avec Monsieur Montemayor."
$ws.Range("B4").Value = "Creation du Git + IceScrum"
$ws.Range("E4").Value = "Prise en main de icescrum."
$ws.Range("B5").Value = "Realisation des sprints sur papier"
$ws.Range("B6").Value = "Sprints + tasks IceScrum"
$ws.Range("E6").Value = "Planification initiale."
$ws.Range("B7").Value = "Creation des premiers documents"
$ws.Range("E7").Value = "Plus debut de la mise en page de la Doc."
$ws.Range("B8").Value = "Preparation du mail d'envoie"
$ws.Range("B9").Value = "MCD-MLD papier"
$ws.Range("B10").Value = "MCD numerique"
$ws.Range("B11").Value = "MLD numerique"
$ws.Range("B12").Value = "Preparation du template / maquette visuel"
$ws.Range("E12").Value = "Utilisation d'un logiciel nommé nicepage pour crée un template. Il me faut demander a mon chef de projet si il considere ok pour une maquette visuel."
$ws.Range("B13").Value = "Arborescence du site"
$ws.Range("B14").Value = "Documentation"
$ws.Range("E14").Value = "Mise en page + premier point de la docs ecrit."
$ws.Range("B15").Value = "Documentation de mon mcd-mld"
$ws.Range("E15").Value = "Ecriture des points dans la docs de chaque tables des mcd-mld."
$ws.Range("B16").Value = "Ecriture du mail d'envoie des PDF"
$ws.Range("E16").Value = "Comme demandé par Monsieur Montemayor, un envoie au format PDF de la documentation et du journal de travail. J'ai vite remis en forme un bug sur la doc."
$ws.Range("B17").Value = "Création de la maquettes visuel de tout le site"
$ws.Range("E17").Value = "Cela m'a pris pas mal de temps, J'ai fais les maquettes pour l'acceuil, Mon calendrier vue mois et semaine, modifier/ajouter évènement, contact, login, register, profil"
$ws.Range("B18").Value = "Discussion avec le chef de projet"
$ws.Range("E18").Value = "A propos de la pertinance d'un gantt dans un projet en mode Agile. Resultat: Pas besoin j'ai fais juste."
$ws.Range("B19").Value = "Ajout des maquettes et commentaires dans la documentation"
$ws.Range("B20").Value = "Documentation"
$ws.Range("E20").Value = "Debut de usercase + test + autres points de la doc de l'analyse et conception."
$ws.Range("B21").Value = "Derniere retouche sur les maquettes graphique"
$ws.Range("E21").Value = "2-3 mots a passer d'anglais a français + bug d'affichage."
$ws.Range("B22").Value = "Recherche a propos du code du calendrier"
$ws.Range("E22").Value = "Etant donné que j'ai fini le sprint 1, je me documente déjà sur comment va se passer la suite histoire de me donné une idée de combien de temps cela va prendre. J'ai déjà reussi a trouver un code source a étudier: https://codes-sources.commentcamarche.net/source/50541-calendrier-agenda-tres-simple-gerer-les-jours-feries-et-les-jours-speciaux"
$ws.Range("B23").Value = "Ajout d'un ReadMe sur le Git"
$ws.Range("E23").Value = "Ajout + ecritude de celui-ci + redecouverte du MarkDown."
$ws.Range("B24").Value = "Retrospective Sprint 1 sur le Git"
$ws.Range("E24").Value = "en + passage du git en public, sinon la creation d'un git n'est pas possible."
$ws.Range("B25").Value = "Debut du touchage de code avec le calendrier"
$ws.Range("E25").Value = "En attendant la fin du sprint qui est la fin de cette journée. Couleur du calendrier, test des fonctionalité, design, etc. Pas de code majeur juste de la mise en page."

# --- New sprint-2 rows 26-32: Tache, Date, Temps, Description supplementaire ---
$ws.Range("B26").Value = "Ecriture des test et taches du sprint 2"
$ws.Range("C26").Value = 44326
$ws.Range("D26").Value = 60
$ws.Range("B27").Value = "Sprint 1 review"
$ws.Range("C27").Value = 44326
$ws.Range("D27").Value = 60
$ws.Range("B28").Value = "Création de la base de donnée"
$ws.Range("C28").Value = 44326
$ws.Range("D28").Value = 30
$ws.Range("E28").Value = "Base de donnée + utilisateur."
$ws.Range("B29").Value = "Modification dans la doc"
$ws.Range("C29").Value = 44326
$ws.Range("D29").Value = 60
$ws.Range("E29").Value = "Les modifications sont en rapport a ce d'ont nous avons parler pendant la sprint review."
$ws.Range("B30").Value = "Architecture MVC"
$ws.Range("C30").Value = 44326
$ws.Range("D30").Value = 120
$ws.Range("E30").Value = "Adaptation du Template en MVC + Redirection pour une navigation basique sur le site."
$ws.Range("B31").Value = "Diagramme de code pour ajouter un évènement"
$ws.Range("C31").Value = 44326
$ws.Range("D31").Value = 60
$ws.Range("B32").Value = "Documentation"
$ws.Range("C32").Value = 44326
$ws.Range("D32").Value = 30
$ws.Range("E32").Value = "Charte graphique + 2-3 autres points"

# --- Row 29 style fix: B29 drops wrap text (s18 -> s7), E29 drops vertical-center while keeping wrap (s22 -> s4) ---
$ws.Range("B29").WrapText = $false
$ws.Range("E29").VerticalAlignment = -4107

# --- Row heights for the newly populated multi-line rows ---
$ws.Rows.Item(29).RowHeight = 30
$ws.Rows.Item(30).RowHeight = 30
$ws.Rows.Item(31).RowHeight = 30

# --- Scroll / selection to mirror the author leaving the cursor at E33 after the new entries ---
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E33").Select()
